$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Plain numeric values that need no special styling (break/hours columns)
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 0.25
$ws.Range("E14").Value = 0.3
$ws.Range("E17").Value = 0.3
$ws.Range("E19").Value = 0.25
$ws.Range("E20").Value = 0.4
$ws.Range("D21").Value = 0.3
$ws.Range("E22").Value = 0.2
$ws.Range("E23").Value = 0.4

# ---------------------------------------------------------------------------
# 2) Text (shared-string) values -- written in the same order the original
#    author entered them so the sharedStrings table comes out in the same
#    sequence.
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "10/11/2020  Project 2"
$ws.Range("F19").Value = "setting up project. "
$ws.Range("F11").Value = "writing tutorial"
$ws.Range("F20").Value = "Coding player movement and camera follow script"
$ws.Range("F22").Value = "writing and adding shooting script and adding emitter"
$ws.Range("F23").Value = "making bullet and enemy disapear on collision"
$ws.Range("F14").Value = "writing tutorial"
$ws.Range("F17").Value = "writing tutorial"

# ---------------------------------------------------------------------------
# 3) Date values (column A) for new / extended rows
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = 44124
$ws.Range("A11").Value = 44124
$ws.Range("A13").Value = 44131
$ws.Range("A14").Value = 44131
$ws.Range("A16").Value = 44138
$ws.Range("A17").Value = 44138
$ws.Range("A20").Value = 44152
$ws.Range("A21").Value = 44152
$ws.Range("A22").Value = 44152
$ws.Range("A23").Value = 44166

# Apply the existing "date" number format (same as A2) to those cells by
# copying the format only -- this re-uses the workbook's existing style
# instead of inventing a new one.
$ws.Range("A2").Copy()
$ws.Range("A10:A11").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A20:A22").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A23").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Start / End time values (columns B & C) for new rows
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = 0.4826388888888889
$ws.Range("C11").Value = 0.5
$ws.Range("B14").Value = 0.54861111111111105
$ws.Range("C14").Value = 0.56944444444444442
$ws.Range("B17").Value = 0.47222222222222227
$ws.Range("C17").Value = 0.49305555555555558
$ws.Range("B19").Value = 0.46180555555555558
$ws.Range("C19").Value = 0.47916666666666669
$ws.Range("B20").Value = 0.43055555555555558
$ws.Range("C20").Value = 0.45833333333333331
$ws.Range("B22").Value = 0.47916666666666669
$ws.Range("C22").Value = 0.49305555555555558
$ws.Range("B23").Value = 0.4375
$ws.Range("C23").Value = 0.46527777777777773

# Re-use the existing "time" number format (same as B2:C2).
$ws.Range("B2:C2").Copy()
$ws.Range("B11:C11").PasteSpecial(-4122)
$ws.Range("B2:C2").Copy()
$ws.Range("B14:C14").PasteSpecial(-4122)
$ws.Range("B2:C2").Copy()
$ws.Range("B17:C17").PasteSpecial(-4122)
$ws.Range("B2:C2").Copy()
$ws.Range("B19:C19").PasteSpecial(-4122)
$ws.Range("B2:C2").Copy()
$ws.Range("B20:C20").PasteSpecial(-4122)
$ws.Range("B2:C2").Copy()
$ws.Range("B22:C22").PasteSpecial(-4122)
$ws.Range("B2:C2").Copy()
$ws.Range("B23:C23").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5) A19 holds descriptive text, not a real date, but carries a date-style
#    number format plus a top vertical alignment (for the taller row).
# ---------------------------------------------------------------------------
$ws.Range("A19").NumberFormat = "m/d/yyyy"
$ws.Range("A19").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# 6) Header row formatting -- bold + underline
# ---------------------------------------------------------------------------
$ws.Range("B1:F1").Font.Bold = $true
$ws.Range("B1:F1").Font.Underline = $true

# ---------------------------------------------------------------------------
# 7) Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.140625
$ws.Columns.Item(6).ColumnWidth = 81.140625

# ---------------------------------------------------------------------------
# 8) Row height for row 19 (wraps a long description)
# ---------------------------------------------------------------------------
$ws.Rows.Item(19).RowHeight = 58.5

# ---------------------------------------------------------------------------
# 9) Page setup
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 10) Leave the selection where the next blank entry would go.
# ---------------------------------------------------------------------------
$ws.Range("F26").Select()
